$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2918.1428
$ws.Range("I40").Value = 2725
$ws.Range("J40").Value = 2995.4
$ws.Range("K40").Value = 2725
$ws.Range("L40").Value = 2995.4
$ws.Range("M40").Value = -2550
$ws.Range("N40").Value = -3345.4

$ws.Range("H53").Value = 6733
$ws.Range("I53").Value = 6733
$ws.Range("K53").Value = 6733
$ws.Range("M53").Value = -6096

$ws.Range("H62").Value = 23164.25
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 23164.25
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H132").Value = 1022.5814
$ws.Range("I132").Value = 893.08826
$ws.Range("K132").Value = 2679.26478
$ws.Range("M132").Value = -149.26478

$ws.Range("H137").Value = 60377.41
$ws.Range("I137").Value = 985.7143
$ws.Range("J137").Value = 101951.6
$ws.Range("K137").Value = 2957.1429
$ws.Range("L137").Value = 305854.8
$ws.Range("M137").Value = -407.1428999999998
$ws.Range("N137").Value = -310954.8

$ws.Range("H138").Value = 1525.15
$ws.Range("I138").Value = 1193.0625
$ws.Range("J138").Value = 1681.4265
$ws.Range("K138").Value = 3579.1875
$ws.Range("L138").Value = 5044.279500000001
$ws.Range("M138").Value = 1560.8125
$ws.Range("N138").Value = -15324.2795

$ws.Range("H139").Value = 73176.55499999999
$ws.Range("J139").Value = 73176.55499999999
$ws.Range("L139").Value = 73176.55499999999
$ws.Range("N139").Value = -83456.55499999999

$ws.Range("H141").Value = 2002539.5
$ws.Range("I141").Value = 2547141.2
$ws.Range("J141").Value = 5666.3335
$ws.Range("K141").Value = 7641423.600000001
$ws.Range("L141").Value = 16999.0005
$ws.Range("M141").Value = -7636243.600000001
$ws.Range("N141").Value = -27359.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6725.489
$ws.Range("I32").Value = 4332.873
$ws.Range("J32").Value = 15666.315
$ws.Range("K32").Value = 4332.873
$ws.Range("L32").Value = 15666.315
$ws.Range("M32").Value = -4045.873
$ws.Range("N32").Value = -16240.315

$ws.Range("H61").Value = 28690.334
$ws.Range("I61").Value = 37225.227
$ws.Range("K61").Value = 37225.227
$ws.Range("M61").Value = -37013.227

$ws.Range("H110").Value = 1067.7059
$ws.Range("I110").Value = 701
$ws.Range("K110").Value = 701
$ws.Range("M110").Value = 1344

$ws.Range("H136").Value = 28690.334
$ws.Range("I136").Value = 37225.227
$ws.Range("K136").Value = 111675.681
$ws.Range("M136").Value = -109125.681

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5558803
$ws.Range("J20").Value = 4160.75
$ws.Range("L20").Value = 4160.75
$ws.Range("N20").Value = -4654.75

$ws.Range("H86").Value = 334777.34
$ws.Range("I86").Value = 888.3333
$ws.Range("J86").Value = 668666.3
$ws.Range("K86").Value = 888.3333
$ws.Range("L86").Value = 668666.3
$ws.Range("M86").Value = 234.6667
$ws.Range("N86").Value = -670912.3

$ws.Range("H89").Value = 334777.34
$ws.Range("I89").Value = 888.3333
$ws.Range("J89").Value = 668666.3
$ws.Range("K89").Value = 4441.6665
$ws.Range("L89").Value = 3343331.5
$ws.Range("M89").Value = 1174.3335
$ws.Range("N89").Value = -3354563.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 212.375
$ws.Range("I7").Value = 235.71428
$ws.Range("K7").Value = 235.71428
$ws.Range("M7").Value = -122.71428

$ws.Range("H31").Value = 2570.8
$ws.Range("I31").Value = 1554.3846
$ws.Range("J31").Value = 4458.4287
$ws.Range("K31").Value = 1554.3846
$ws.Range("L31").Value = 4458.4287
$ws.Range("M31").Value = -1259.3846
$ws.Range("N31").Value = -5048.4287

$ws.Range("H34").Value = 2570.8
$ws.Range("I34").Value = 1554.3846
$ws.Range("J34").Value = 4458.4287
$ws.Range("K34").Value = 1554.3846
$ws.Range("L34").Value = 4458.4287
$ws.Range("M34").Value = -1352.3846
$ws.Range("N34").Value = -4862.4287

$ws.Range("H107").Value = 1164.6786
$ws.Range("I107").Value = 1062.4286
$ws.Range("J107").Value = 1471.4286
$ws.Range("K107").Value = 1062.4286
$ws.Range("L107").Value = 1471.4286
$ws.Range("M107").Value = 857.5714
$ws.Range("N107").Value = -5311.4286

$ws.Range("H122").Value = 1613.5
$ws.Range("I122").Value = 1613.5
$ws.Range("K122").Value = 4840.5
$ws.Range("M122").Value = -2390.5

$ws.Range("H132").Value = 1426.1562
$ws.Range("I132").Value = 808.8889
$ws.Range("K132").Value = 2426.6667
$ws.Range("M132").Value = 103.3332999999998

$ws.Range("H133").Value = 30350
$ws.Range("J133").Value = 30350
$ws.Range("L133").Value = 30350
$ws.Range("N133").Value = -35410

$ws.Range("H135").Value = 34694
$ws.Range("J135").Value = 34694
$ws.Range("L135").Value = 34694
$ws.Range("N135").Value = -44834

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 194.85715
$ws.Range("I33").Value = 112.375
$ws.Range("J33").Value = 304.83334
$ws.Range("K33").Value = 674.25
$ws.Range("L33").Value = 1829.00004
$ws.Range("M33").Value = -391.25
$ws.Range("N33").Value = -2395.00004

$ws.Range("H105").Value = 6610.9443
$ws.Range("J105").Value = 6610.9443
$ws.Range("L105").Value = 19832.8329
$ws.Range("N105").Value = -25074.8329

$ws.Range("H126").Value = 3750

$ws.Range("H131").Value = 14939.375
$ws.Range("J131").Value = 17021.877
$ws.Range("L131").Value = 51065.631
$ws.Range("N131").Value = -61145.631

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws.Range("H62").Value = 33998.5
$ws.Range("J62").Value = 33998.5
$ws.Range("L62").Value = 33998.5
$ws.Range("N62").Value = -35370.5

$ws.Range("H65").Value = 33998.5
$ws.Range("J65").Value = 33998.5
$ws.Range("L65").Value = 101995.5
$ws.Range("N65").Value = -108859.5

$ws.Range("H80").Value = 2848.5
$ws.Range("I80").Value = 2631.3333
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 2631.3333
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = -1633.3333
$ws.Range("N80").Value = -5496

$ws.Range("H83").Value = 2848.5
$ws.Range("I83").Value = 2631.3333
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 13156.6665
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = -8164.666499999999
$ws.Range("N83").Value = -27484

$ws.Range("H102").Value = 2749.3333
$ws.Range("I102").Value = 2749.3333
$ws.Range("K102").Value = 2749.3333
$ws.Range("M102").Value = -1127.3333

$ws.Range("H122").Value = 3975
$ws.Range("I122").Value = 4800
$ws.Range("K122").Value = 14400
$ws.Range("M122").Value = -11950

$ws.Range("H126").Value = 2461772
$ws.Range("I126").Value = 2927076.5
$ws.Range("K126").Value = 8781229.5
$ws.Range("M126").Value = -8778759.5

$ws.Range("H136").Value = 18853.637
$ws.Range("J136").Value = 18853.637
$ws.Range("L136").Value = 56560.91099999999
$ws.Range("N136").Value = -61660.91099999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1166.3529
$ws.Range("I22").Value = 805.875
$ws.Range("J22").Value = 1486.7778
$ws.Range("K22").Value = 805.875
$ws.Range("L22").Value = 1486.7778
$ws.Range("M22").Value = -510.875
$ws.Range("N22").Value = -2076.7778

$ws.Range("H27").Value = 1166.3529
$ws.Range("I27").Value = 805.875
$ws.Range("J27").Value = 1486.7778
$ws.Range("K27").Value = 805.875
$ws.Range("L27").Value = 1486.7778
$ws.Range("M27").Value = -698.875
$ws.Range("N27").Value = -1700.7778

$ws.Range("H93").Value = 1220.8077
$ws.Range("I93").Value = 714.4737
$ws.Range("K93").Value = 714.4737
$ws.Range("M93").Value = 533.5263

$ws.Range("H136").Value = 3755.7778
$ws.Range("I136").Value = 2929
$ws.Range("K136").Value = 8787
$ws.Range("M136").Value = -6237

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 30808
$ws.Range("J16").Value = 30808
$ws.Range("L16").Value = 30808
$ws.Range("N16").Value = -31392

$ws.Range("H40").Value = 57028
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H100").Value = 1189
$ws.Range("I100").Value = 900.5
$ws.Range("K100").Value = 1801
$ws.Range("M100").Value = -1260

$ws.Range("H107").Value = 625.7059
$ws.Range("I107").Value = 539.8125
$ws.Range("K107").Value = 1619.4375
$ws.Range("M107").Value = 300.5625

$ws.Range("H126").Value = 1347.4849
$ws.Range("I126").Value = 1097.92
$ws.Range("K126").Value = 3293.76
$ws.Range("M126").Value = -823.7600000000002

$ws.Range("H132").Value = 1248.4154
$ws.Range("I132").Value = 1185.9231
$ws.Range("K132").Value = 3557.7693
$ws.Range("M132").Value = -1027.7693
